# The commit swaps the presentation's design/theme: the slide master's
# theme ("Integral" colours) is replaced with the built-in "Office Theme"
# colour palette (the presentation's previously-cached alternate theme).
#
# PowerPoint exposes the active theme's 12 scheme colours via
# Master.Theme.ThemeColorScheme.Colors(index).RGB (msoThemeColorSchemeIndex
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Re-pointing every
# slot to the Office Theme RGB values reproduces the colour swap the author
# made from PowerPoint's Design gallery.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$officeThemeColors = @(
    0,         # 1  dk1       000000
    16777215,  # 2  lt1       FFFFFF
    6968388,   # 3  dk2       44546A
    15132391,  # 4  lt2       E7E6E6
    13998939,  # 5  accent1   5B9BD5
    3243501,   # 6  accent2   ED7D31
    10855845,  # 7  accent3   A5A5A5
    49407,     # 8  accent4   FFC000
    12874308,  # 9  accent5   4472C4
    4697456,   # 10 accent6   70AD47
    12673797,  # 11 hlink     0563C1
    7491477    # 12 folHlink  954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
